# fix(employments): displaying employments details data only if provided
#
# The "Employment" sheet shows, per actor type, five employment-category
# columns (Temporary Male/Female, Permanent Unskilled Male/Female,
# Permanent Skilled Male/Female) each paired with a "value" column
# (B, D, F, H, J, L) and a "unit" column (C, E, G, I, K, M).
#
# Only the "Temporary Male" (B/C) and "Permanent Unskilled Male" (F/G)
# figures were actually provided in the source data; the other value
# columns (D, H, J, L) were left completely blank instead of showing a
# placeholder. This fills those provided-but-empty columns with the same
# "-" placeholder already used elsewhere in the workbook, so the sheet
# only *displays* employment details data for the categories that have a
# value (the truly missing ones render "-" instead of an empty cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employment")

$firstRow = 2
$lastRow = 12
$emptyValueColumns = @("D", "H", "J", "L")

for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($col in $emptyValueColumns) {
        $ws.Range("$col$r").Value = "-"
    }
}

# Give the "Actor types" sheet an explicit width for column A (it was
# relying on the default width before).
$actorTypes = $wb.Worksheets.Item("Actor types")
$actorTypes.Columns.Item(1).ColumnWidth = 13.5

# Make "Employment" the active sheet/tab, with L2:L12 (the newly-filled
# placeholder column) selected.
$ws.Activate()
[void]$ws.Range("L2:L12").Select()
